# Update IG attribution values for rows 2 and 3 (relative direction fix)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 0.005130496560918854
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = -0
$ws.Cells.Item(2, 4).Value = -0.02392885496758202
$ws.Cells.Item(2, 5).Value = -0.01693560061416057
$ws.Cells.Item(2, 6).Value = -0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = -0
$ws.Cells.Item(2, 10).Value = 0.02471286464988766
$ws.Cells.Item(2, 11).Value = -0
$ws.Cells.Item(2, 12).Value = -0
$ws.Cells.Item(2, 13).Value = -0.004547601632370699
$ws.Cells.Item(2, 14).Value = -0.02502063099885841
$ws.Cells.Item(2, 17).Value = -0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = -0.02766699180611887
$ws.Cells.Item(2, 20).Value = -0
$ws.Cells.Item(2, 21).Value = -0
$ws.Cells.Item(2, 22).Value = 0.05419437097972913
$ws.Cells.Item(2, 23).Value = 0.02599247979354504
$ws.Cells.Item(2, 27).Value = 0
$ws.Cells.Item(2, 28).Value = 0.004743365091757061
$ws.Cells.Item(2, 29).Value = -0
$ws.Cells.Item(2, 30).Value = -0
$ws.Cells.Item(2, 31).Value = -0.02707591191183536
$ws.Cells.Item(2, 32).Value = 0.01709765687544249
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 37).Value = 0.01784621566429674
$ws.Cells.Item(2, 38).Value = 0
$ws.Cells.Item(2, 39).Value = 0
$ws.Cells.Item(2, 40).Value = 0.002509656631344049
$ws.Cells.Item(2, 41).Value = -0.0591894262847708
$ws.Cells.Item(2, 44).Value = 0
$ws.Cells.Item(2, 46).Value = 0.02245941889391236
$ws.Cells.Item(2, 47).Value = 0
$ws.Cells.Item(2, 48).Value = 0
$ws.Cells.Item(2, 49).Value = -0.03997207371561685
$ws.Cells.Item(2, 50).Value = -0.0222455408604398
$ws.Cells.Item(2, 51).Value = -0
$ws.Cells.Item(2, 55).Value = 0.02500263180353006
$ws.Cells.Item(2, 56).Value = 0
$ws.Cells.Item(2, 57).Value = -0
$ws.Cells.Item(2, 58).Value = -0.002168975133050314
$ws.Cells.Item(2, 59).Value = -0.01063924078124795
$ws.Cells.Item(2, 61).Value = 0
$ws.Cells.Item(2, 62).Value = -0
$ws.Cells.Item(2, 63).Value = -0
$ws.Cells.Item(2, 64).Value = 0.02326066322100957
$ws.Cells.Item(2, 65).Value = 0
$ws.Cells.Item(2, 66).Value = -0
$ws.Cells.Item(2, 67).Value = 0.03178992770114826
$ws.Cells.Item(2, 68).Value = 0.03817303067736682
$ws.Cells.Item(2, 69).Value = -0
$ws.Cells.Item(2, 71).Value = -0
$ws.Cells.Item(2, 73).Value = -0.001857381766410295
$ws.Cells.Item(2, 74).Value = -0
$ws.Cells.Item(2, 75).Value = -0
$ws.Cells.Item(2, 76).Value = -0.01031704585172784
$ws.Cells.Item(2, 77).Value = -0.02132798982982983
$ws.Cells.Item(2, 78).Value = -0
$ws.Cells.Item(2, 82).Value = -0.007239308633253068
$ws.Cells.Item(2, 83).Value = 0
$ws.Cells.Item(2, 84).Value = 0
$ws.Cells.Item(2, 85).Value = 0.001301096967254109
$ws.Cells.Item(2, 86).Value = 0.01597316963689811
$ws.Cells.Item(2, 87).Value = 0
$ws.Cells.Item(2, 88).Value = -0
$ws.Cells.Item(2, 89).Value = -0
$ws.Cells.Item(2, 90).Value = -0
$ws.Cells.Item(2, 91).Value = -0.004172243429578691
$ws.Cells.Item(2, 92).Value = 0
$ws.Cells.Item(2, 93).Value = -0
$ws.Cells.Item(2, 94).Value = -0.004981208171710785
$ws.Cells.Item(2, 95).Value = -0.01577156610364081
$ws.Cells.Item(2, 96).Value = -0
$ws.Cells.Item(2, 97).Value = -0
$ws.Cells.Item(2, 100).Value = -0.00217257729493647
$ws.Cells.Item(2, 101).Value = -0
$ws.Cells.Item(2, 102).Value = -0
$ws.Cells.Item(2, 103).Value = 0.02182022279621677
$ws.Cells.Item(2, 104).Value = 0.02569762313996425
$ws.Cells.Item(2, 107).Value = 0
$ws.Cells.Item(2, 109).Value = -0.0002172611616025855
$ws.Cells.Item(2, 110).Value = 0
$ws.Cells.Item(2, 111).Value = 0
$ws.Cells.Item(2, 112).Value = 0.08175342688176625
$ws.Cells.Item(2, 113).Value = -0.05464954313494159
$ws.Cells.Item(2, 115).Value = -0
$ws.Cells.Item(2, 118).Value = -0.01734652580522425
$ws.Cells.Item(2, 119).Value = -0
$ws.Cells.Item(2, 120).Value = -0
$ws.Cells.Item(2, 121).Value = 0.0114803559676882
$ws.Cells.Item(2, 122).Value = 0.1013287922281259
$ws.Cells.Item(2, 124).Value = -0
$ws.Cells.Item(2, 125).Value = -0
$ws.Cells.Item(2, 127).Value = -0.01179741443645917
$ws.Cells.Item(2, 128).Value = 0
$ws.Cells.Item(2, 129).Value = 0
$ws.Cells.Item(2, 130).Value = -0.01433721291011771
$ws.Cells.Item(2, 131).Value = 0.02473195073079327
$ws.Cells.Item(2, 132).Value = -0
$ws.Cells.Item(2, 133).Value = -0
$ws.Cells.Item(2, 134).Value = 0
$ws.Cells.Item(2, 136).Value = -0.01408018348142348
$ws.Cells.Item(2, 137).Value = 0
$ws.Cells.Item(2, 138).Value = -0
$ws.Cells.Item(2, 139).Value = 0.0365070358046652
$ws.Cells.Item(2, 140).Value = 0.07043197269217763
$ws.Cells.Item(2, 141).Value = 0
$ws.Cells.Item(2, 145).Value = -0.01767988419955132
$ws.Cells.Item(2, 146).Value = -0
$ws.Cells.Item(2, 147).Value = 0
$ws.Cells.Item(2, 148).Value = 0.01408851028855623
$ws.Cells.Item(2, 149).Value = 0.04907402152391164
$ws.Cells.Item(2, 150).Value = 0
$ws.Cells.Item(2, 153).Value = -0
$ws.Cells.Item(2, 154).Value = -0.00326685600030628
$ws.Cells.Item(2, 155).Value = 0
$ws.Cells.Item(2, 156).Value = 0
$ws.Cells.Item(2, 157).Value = 0.01708878651669269
$ws.Cells.Item(2, 158).Value = 0.02962248912172181
$ws.Cells.Item(2, 161).Value = -0
$ws.Cells.Item(2, 162).Value = -0
$ws.Cells.Item(2, 163).Value = 0.01353753197573294
$ws.Cells.Item(2, 164).Value = 0
$ws.Cells.Item(2, 165).Value = -0
$ws.Cells.Item(2, 166).Value = 0.0473094570095499
$ws.Cells.Item(2, 167).Value = 0.07124980721699638
$ws.Cells.Item(2, 168).Value = 0
$ws.Cells.Item(2, 169).Value = -0
$ws.Cells.Item(2, 170).Value = -0
$ws.Cells.Item(2, 172).Value = -0.01169197431578245
$ws.Cells.Item(2, 173).Value = 0
$ws.Cells.Item(2, 174).Value = 0
$ws.Cells.Item(2, 175).Value = -0.009271809757468087
$ws.Cells.Item(2, 176).Value = 0.08485382841411168
$ws.Cells.Item(2, 177).Value = 0
$ws.Cells.Item(2, 181).Value = -0.0007367975867645503
$ws.Cells.Item(2, 182).Value = 0
$ws.Cells.Item(2, 183).Value = -0
$ws.Cells.Item(2, 184).Value = 0.03456136806104576
$ws.Cells.Item(2, 185).Value = 0
$ws.Cells.Item(2, 186).Value = 0
$ws.Cells.Item(2, 187).Value = -0
$ws.Cells.Item(2, 188).Value = 0
$ws.Cells.Item(2, 189).Value = -0

# Row 3
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 3).Value = -0
$ws.Cells.Item(3, 4).Value = -0.01802642255756185
$ws.Cells.Item(3, 5).Value = 0.00137538753950138
$ws.Cells.Item(3, 6).Value = 0.5052551783512984
$ws.Cells.Item(3, 7).Value = -0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = -0
$ws.Cells.Item(3, 12).Value = -0
$ws.Cells.Item(3, 13).Value = -0.0303317485894369
$ws.Cells.Item(3, 14).Value = -0.02436577158335003
$ws.Cells.Item(3, 15).Value = 0.3482606298072715
$ws.Cells.Item(3, 16).Value = -0
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 22).Value = 0.01663914243282382
$ws.Cells.Item(3, 23).Value = 0.07818631825847189
$ws.Cells.Item(3, 24).Value = 0.04737148275777644
$ws.Cells.Item(3, 25).Value = -0
$ws.Cells.Item(3, 26).Value = 0
$ws.Cells.Item(3, 27).Value = 0
$ws.Cells.Item(3, 28).Value = -0
$ws.Cells.Item(3, 31).Value = -0.02801164544941828
$ws.Cells.Item(3, 32).Value = 0.02602505742690114
$ws.Cells.Item(3, 33).Value = 0.02708899519688088
$ws.Cells.Item(3, 36).Value = 0
$ws.Cells.Item(3, 37).Value = -0
$ws.Cells.Item(3, 38).Value = 0
$ws.Cells.Item(3, 40).Value = -0.02702469349752857
$ws.Cells.Item(3, 41).Value = -0.0594813019651976
$ws.Cells.Item(3, 42).Value = 0.04579260220400171
$ws.Cells.Item(3, 44).Value = 0
$ws.Cells.Item(3, 45).Value = 0
$ws.Cells.Item(3, 46).Value = -0
$ws.Cells.Item(3, 47).Value = -0
$ws.Cells.Item(3, 49).Value = -0.02005802770442524
$ws.Cells.Item(3, 50).Value = 0.04886987481839743
$ws.Cells.Item(3, 51).Value = 0.328190777470861
$ws.Cells.Item(3, 54).Value = -0
$ws.Cells.Item(3, 55).Value = 0
$ws.Cells.Item(3, 56).Value = 0
$ws.Cells.Item(3, 58).Value = -0.01604642426108615
$ws.Cells.Item(3, 59).Value = -0.00488960269155806
$ws.Cells.Item(3, 60).Value = 0.08156072267757562
$ws.Cells.Item(3, 62).Value = 0
$ws.Cells.Item(3, 63).Value = 0
$ws.Cells.Item(3, 64).Value = -0
$ws.Cells.Item(3, 65).Value = 0
$ws.Cells.Item(3, 66).Value = 0
$ws.Cells.Item(3, 67).Value = 0.01078948871334732
$ws.Cells.Item(3, 68).Value = 0.02224153603180841
$ws.Cells.Item(3, 69).Value = 0.04228142164201944
$ws.Cells.Item(3, 70).Value = -0
$ws.Cells.Item(3, 71).Value = -0
$ws.Cells.Item(3, 72).Value = -0
$ws.Cells.Item(3, 73).Value = -0
$ws.Cells.Item(3, 74).Value = -0
$ws.Cells.Item(3, 76).Value = -0.02173066233972291
$ws.Cells.Item(3, 77).Value = 0.03998230953093256
$ws.Cells.Item(3, 78).Value = 0.1320242281631021
$ws.Cells.Item(3, 79).Value = -0
$ws.Cells.Item(3, 80).Value = 0
$ws.Cells.Item(3, 81).Value = -0
$ws.Cells.Item(3, 82).Value = -0
$ws.Cells.Item(3, 84).Value = 0
$ws.Cells.Item(3, 85).Value = 0.01894025474399924
$ws.Cells.Item(3, 86).Value = -0.0008720681587888977
$ws.Cells.Item(3, 87).Value = -0.0970636524279626
$ws.Cells.Item(3, 88).Value = 0
$ws.Cells.Item(3, 90).Value = 0
$ws.Cells.Item(3, 91).Value = -0
$ws.Cells.Item(3, 93).Value = 0
$ws.Cells.Item(3, 94).Value = -0.003492794481246962
$ws.Cells.Item(3, 95).Value = -0.02115982449754519
$ws.Cells.Item(3, 96).Value = -0.02373757270871685
$ws.Cells.Item(3, 97).Value = 0
$ws.Cells.Item(3, 98).Value = -0
$ws.Cells.Item(3, 99).Value = 0
$ws.Cells.Item(3, 100).Value = 0
$ws.Cells.Item(3, 103).Value = 0.0244168598420817
$ws.Cells.Item(3, 104).Value = 0.004767563638659294
$ws.Cells.Item(3, 105).Value = -0.08992352107104465
$ws.Cells.Item(3, 108).Value = -0
$ws.Cells.Item(3, 109).Value = -0
$ws.Cells.Item(3, 111).Value = 0
$ws.Cells.Item(3, 112).Value = -0.01984978673041322
$ws.Cells.Item(3, 113).Value = -0.02627699406989251
$ws.Cells.Item(3, 114).Value = -0.02077541640713792
$ws.Cells.Item(3, 115).Value = 0
$ws.Cells.Item(3, 117).Value = -0
$ws.Cells.Item(3, 118).Value = -0
$ws.Cells.Item(3, 119).Value = 0
$ws.Cells.Item(3, 121).Value = -0.01503995111913278
$ws.Cells.Item(3, 122).Value = 0.07076097943013691
$ws.Cells.Item(3, 123).Value = 0.01380463233941497
$ws.Cells.Item(3, 125).Value = 0
$ws.Cells.Item(3, 126).Value = 0
$ws.Cells.Item(3, 127).Value = -0
$ws.Cells.Item(3, 128).Value = 0
$ws.Cells.Item(3, 129).Value = 0
$ws.Cells.Item(3, 130).Value = -0.03230563932452513
$ws.Cells.Item(3, 131).Value = -0.006141079327499638
$ws.Cells.Item(3, 132).Value = -0.051662033000376
$ws.Cells.Item(3, 134).Value = -0
$ws.Cells.Item(3, 135).Value = -0
$ws.Cells.Item(3, 136).Value = 0
$ws.Cells.Item(3, 137).Value = -0
$ws.Cells.Item(3, 138).Value = 0
$ws.Cells.Item(3, 139).Value = 0.0585366265569453
$ws.Cells.Item(3, 140).Value = 0.0270182361087934
$ws.Cells.Item(3, 141).Value = -0.1086868860668046
$ws.Cells.Item(3, 144).Value = -0
$ws.Cells.Item(3, 145).Value = -0
$ws.Cells.Item(3, 148).Value = 0.04672893403672061
$ws.Cells.Item(3, 149).Value = 0.03172565915530302
$ws.Cells.Item(3, 150).Value = -0.07087105212491739
$ws.Cells.Item(3, 151).Value = 0
$ws.Cells.Item(3, 153).Value = 0
$ws.Cells.Item(3, 154).Value = -0
$ws.Cells.Item(3, 155).Value = 0
$ws.Cells.Item(3, 156).Value = 0
$ws.Cells.Item(3, 157).Value = 0.01557261574551646
$ws.Cells.Item(3, 158).Value = -0.01728727679270572
$ws.Cells.Item(3, 159).Value = -0.08674334406668699
$ws.Cells.Item(3, 162).Value = 0
$ws.Cells.Item(3, 163).Value = -0
$ws.Cells.Item(3, 164).Value = -0
$ws.Cells.Item(3, 166).Value = -0.009839159865455702
$ws.Cells.Item(3, 167).Value = -0.02811539912733575
$ws.Cells.Item(3, 168).Value = 0.0003142644878452895
$ws.Cells.Item(3, 169).Value = 0
$ws.Cells.Item(3, 170).Value = -0
$ws.Cells.Item(3, 171).Value = 0
$ws.Cells.Item(3, 172).Value = 0
$ws.Cells.Item(3, 173).Value = 0
$ws.Cells.Item(3, 174).Value = 0
$ws.Cells.Item(3, 175).Value = -0.01925484836449415
$ws.Cells.Item(3, 176).Value = -0.1268052952161465
$ws.Cells.Item(3, 177).Value = -0.03847016567433773
$ws.Cells.Item(3, 178).Value = 0
$ws.Cells.Item(3, 179).Value = 0
$ws.Cells.Item(3, 180).Value = -0
$ws.Cells.Item(3, 181).Value = -0
$ws.Cells.Item(3, 182).Value = -0
$ws.Cells.Item(3, 184).Value = 0.02741854847410577
$ws.Cells.Item(3, 185).Value = -0
$ws.Cells.Item(3, 186).Value = 0.03229545231933784
$ws.Cells.Item(3, 189).Value = -0
